$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 44, pushing the existing rows 44-45 down to 46-47
$ws.Rows("44:45").Insert()

# Row 44 - new data
$ws.Range("A44").Value = 3
$ws.Range("B44").Value = "Femacal de La Calera"
$ws.Range("C44").Value = "Coquimbo"
$ws.Range("D44").Value = 45166
$ws.Range("D44").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E44").Value = 5
$ws.Range("F44").Value = 100112043
$ws.Range("G44").Value = "Pepino dulce"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 68
$ws.Range("K44").Value = 22000
$ws.Range("L44").Value = 22000
$ws.Range("M44").Value = 22000
$ws.Range("N44").Value = "`$/caja 15 kilos"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 1467
$ws.Range("Q44").Value = 15
$ws.Range("R44").Value = "Hortaliza"

# Row 45 - new data
$ws.Range("A45").Value = 3
$ws.Range("B45").Value = "Femacal de La Calera"
$ws.Range("C45").Value = "Coquimbo"
$ws.Range("D45").Value = 45166
$ws.Range("D45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E45").Value = 5
$ws.Range("F45").Value = 100112043
$ws.Range("G45").Value = "Pepino dulce"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Segunda"
$ws.Range("J45").Value = 50
$ws.Range("K45").Value = 16000
$ws.Range("L45").Value = 16000
$ws.Range("M45").Value = 16000
$ws.Range("N45").Value = "`$/caja 15 kilos"
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 1067
$ws.Range("Q45").Value = 15
$ws.Range("R45").Value = "Hortaliza"
